$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '74.855.60'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '2.813.55'
$ws.Range("E3").Value = '  +7.02%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '596.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +3.21%  '
$ws.Range("E9").Value = '  -5.64%  '
$ws.Range("D10").Value = '2.812.63'
$ws.Range("E10").Value = '  +6.99%  '
$ws.Range("E12").Value = '  +1.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").Value = '3.327.12'
$ws.Range("E14").Value = '  +7.20%  '
$ws.Range("D15").Value = '74.820.45'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("D18").Value = '2.812.04'
$ws.Range("E18").Value = '  +7.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").Value = '2.958.76'
$ws.Range("E27").Value = '  +7.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.68%  '
$ws.Range("E30").Value = '  +7.01%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '513.37'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("E33").Value = '  -1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.86'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("E35").Value = '  +1.89%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.08%  '
$ws.Range("E39").Value = '  -3.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '182.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.96%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.60%  '
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0859'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.68%  '
$ws.Range("E50").Value = '  +6.70%  '
$ws.Range("E51").Value = '  +2.14%  '
